# Update the price-list date in cell A1 of the first sheet (Hoja1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = 45405
